# Fruta / hortaliza, semanal
# Insert a new weekly record as row 241 in the "Vega Modelo de Temuco - Piña" sheet,
# pushing the existing rows 241:267 down to 242:268.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 241 (shifts rows 241-267 down to 242-268)
$ws.Rows.Item(241).Insert()

# Populate the new row 241 with this week's data
$ws.Cells.Item(241, 1).Value() = 10
$ws.Cells.Item(241, 2).Value() = "Vega Modelo de Temuco"
$ws.Cells.Item(241, 3).Value() = "La Araucanía"
$ws.Cells.Item(241, 4).Value() = 44449
$ws.Cells.Item(241, 5).Value() = 9
$ws.Cells.Item(241, 6).Value() = "Fruta"
$ws.Cells.Item(241, 7).Value() = 100108
$ws.Cells.Item(241, 8).Value() = "Tropicales y subtropicales"
$ws.Cells.Item(241, 9).Value() = 100108005
$ws.Cells.Item(241, 10).Value() = "Piña"
$ws.Cells.Item(241, 11).Value() = "Caramelo"
$ws.Cells.Item(241, 12).Value() = "Primera"
$ws.Cells.Item(241, 13).Value() = 50
$ws.Cells.Item(241, 14).Value() = 20000
$ws.Cells.Item(241, 15).Value() = 20000
$ws.Cells.Item(241, 16).Value() = 20000
$ws.Cells.Item(241, 17).Value() = "$/caja 12 unidades"
$ws.Cells.Item(241, 18).Value() = "Ecuador"
$ws.Cells.Item(241, 19).Value() = 1667
$ws.Cells.Item(241, 20).Value() = 12
